$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-point the four allergy/hypersensitivity-related rows (Atopic
#    dermatitis, Allergic rhinitis, Angioedema, Anaphylaxis) at the new
#    grouped "Hypersensitivity" phenotype id, highlight their phenotypeName
#    cell (column D) in yellow, and flag them for deprecation (column I = Y).
# ---------------------------------------------------------------------------
$rowsToGroup = @(26, 70, 123, 171)
foreach ($r in $rowsToGroup) {
    $ws.Cells.Item($r, 2).Value = 43021226000
    $ws.Cells.Item($r, 4).Interior.Color = 65535
    $ws.Cells.Item($r, 9).Value = "Y"
}

# Row 180 (Cancer) previously carried the yellow highlight; it is no longer
# part of a just-grouped set, so the highlight is cleared (explicit "No
# Fill", matching the other already-cleared rows' style in the workbook).
$ws.Range("D6").Copy()
$ws.Range("D180").PasteSpecial(-4122)
$ws.Range("D180").Value2 = $ws.Range("D180").Value2

# ---------------------------------------------------------------------------
# 2. Insert the new "Hypersensitivity" row right before the old row 265.
# ---------------------------------------------------------------------------
$ws.Rows.Item(265).Insert()
$ws.Cells.Item(265, 2).Value = 43021226000
$ws.Cells.Item(265, 3).Value = 43021226000
$ws.Cells.Item(265, 4).Value = 43021226000
$ws.Cells.Item(265, 5).Value = 43021226
$ws.Cells.Item(265, 6).Value = "Overview: Presentation: Assessment: Plan: Prognosis:"
$ws.Cells.Item(265, 8).Value = "To Do"

# ---------------------------------------------------------------------------
# 3. Refresh the AutoFilter / defined name / selection so they cover the new
#    data extent (A1:I269 instead of A1:I268).
# ---------------------------------------------------------------------------
if ($ws.AutoFilterMode) {
    $ws.AutoFilter.Range.AutoFilter()
}
$ws.Range("A1:I269").AutoFilter()

$sortKey = $ws.Range("C1:C269")
$ws.Range("A1:I269").Sort($sortKey, 1, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)

$dn = $wb.Names.Item("PhenotypeDescription!_FilterDatabase")
$dn.RefersTo = "=PhenotypeDescription!`$A`$1:`$I`$269"

$ws.Range("D2").Select()
